$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "last updated" timestamp string (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Septiembre de 2020 a las 18:25"

# ---------------------------------------------------------------------------
# 2) Update updated case counts for the countries whose figures changed.
#    (Row numbers correspond to the already-sorted-by-total-cases layout;
#     country identity/order for these specific rows is unaffected.)
# ---------------------------------------------------------------------------

# Estados Unidos (row 4)
$ws.Range("B4").Value = 7253443
$ws.Range("C4").Value = 9259
$ws.Range("D4").Value = 4484953
$ws.Range("E4").Value = 2559838
$ws.Range("G4").Value = 212
$ws.Range("H4").Value = 208652

# Reino Unido (row 17)
$ws.Range("B17").Value = 429277
$ws.Range("C17").Value = 6042
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = 41971

# Canada (row 29)
$ws.Range("B29").Value = 151589
$ws.Range("C29").Value = 1133
$ws.Range("D29").Value = 130315
$ws.Range("E29").Value = 12011
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = 9263

# Guatemala (row 45)
$ws.Range("B45").Value = 89702
$ws.Range("C45").Value = 824
$ws.Range("D45").Value = 78284
$ws.Range("E45").Value = 8205
$ws.Range("G45").Value = 27
$ws.Range("H45").Value = 3213

# Chequia (row 57)
$ws.Range("B57").Value = 62559
$ws.Range("C57").Value = 1241
$ws.Range("D57").Value = 30891
$ws.Range("E57").Value = 31080
$ws.Range("G57").Value = 7
$ws.Range("H57").Value = 588

# Moldavia (row 63)
$ws.Range("B63").Value = 50534
$ws.Range("C63").Value = 868
$ws.Range("D63").Value = 37440
$ws.Range("E63").Value = 11815
$ws.Range("G63").Value = 15
$ws.Range("H63").Value = 1279

# Grecia (row 87)
$ws.Range("B87").Value = 17228
$ws.Range("C87").Value = 315
$ws.Range("E87").Value = 6863
$ws.Range("G87").Value = 7
$ws.Range("H87").Value = 376

# Luxemburgo (row 109)
$ws.Range("B109").Value = 8311
$ws.Range("C109").Value = 78
$ws.Range("E109").Value = 1211

# ---------------------------------------------------------------------------
# 3) Montenegro's total cases overtook Maldivas and Birmania, so it now
#    ranks ahead of them. Re-write rows 100-102 so that:
#      row 100 -> Montenegro (with its newly updated figures)
#      row 101 -> Maldivas   (figures unchanged)
#      row 102 -> Birmania   (figures unchanged)
# ---------------------------------------------------------------------------
$ws.Range("A100").Value = "Montenegro"
$ws.Range("B100").Value = 10197
$ws.Range("C100").Value = 235
$ws.Range("D100").Value = 6368
$ws.Range("E100").Value = 3671
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 3
$ws.Range("H100").Value = 158

$ws.Range("A101").Value = "Maldivas"
$ws.Range("B101").Value = 10014
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 8673
$ws.Range("E101").Value = 1307
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 34

$ws.Range("A102").Value = "Birmania"
$ws.Range("B102").Value = 9991
$ws.Range("C102").Value = 879
$ws.Range("D102").Value = 2681
$ws.Range("E102").Value = 7112
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 24
$ws.Range("H102").Value = 198

# ---------------------------------------------------------------------------
# 4) Santa Lucia and Timor Oriental are tied on every figure, but swapped
#    their relative order. Swap the country names only (figures identical).
# ---------------------------------------------------------------------------
$ws.Range("A206").Value = "Santa Lucia"
$ws.Range("A207").Value = "Timor Oriental"
